# Week 16 stat log + season sim (Week 17) update
# - Adds a new Rushing row for J.Johnson (inserted after T.Huntley, before T.Williams)
# - Bumps cumulative Rushing totals for L.Murray and D.Freeman
# - Bumps cumulative Receiving totals for several players

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Rushing": shift rows 4..10 down to 5..11 (bottom-up, so we never
# clobber a source row before it's copied), then write the new J.Johnson row
# into row 4 with its own stats.
# ---------------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

for ($r = 10; $r -ge 4; $r--) {
    $dest = $r + 1
    $rushing.Cells.Item($dest, 2).Value = $rushing.Cells.Item($r, 2).Value2
    $rushing.Cells.Item($dest, 3).Value = $rushing.Cells.Item($r, 3).Value2
    $rushing.Cells.Item($dest, 4).Value = $rushing.Cells.Item($r, 4).Value2
    $rushing.Cells.Item($dest, 5).Value = $rushing.Cells.Item($r, 5).Value2
    $rushing.Cells.Item($dest, 6).Value = $rushing.Cells.Item($r, 6).Value2
    # Column A is a plain 0-based row index, independent of the shifted
    # player - renumber it in place so it stays a contiguous sequence.
    $rushing.Cells.Item($dest, 1).Value = $dest - 2
}

# Row 11 is brand-new (sheet previously ended at row 10) - copy the A-column
# number-style formatting down from row 10 so it matches the rest of the
# table instead of picking up a blank default style.
$rushing.Cells.Item(10, 1).Copy()
$rushing.Cells.Item(11, 1).PasteSpecial(-4122)

# New player row: J.Johnson
$rushing.Cells.Item(4, 1).Value = 2
$rushing.Cells.Item(4, 2).Value = "J.Johnson"
$rushing.Cells.Item(4, 3).Value = 1
$rushing.Cells.Item(4, 4).Value = 1
$rushing.Cells.Item(4, 5).Value = 3
$rushing.Cells.Item(4, 6).Value = 1

# Week 16 stat bumps for existing rushers (now one row lower than before)
$rushing.Cells.Item(6, 3).Value = 51
$rushing.Cells.Item(6, 4).Value = 31
$rushing.Cells.Item(6, 5).Value = 9
$rushing.Cells.Item(6, 6).Value = 16

$rushing.Cells.Item(7, 3).Value = 65
$rushing.Cells.Item(7, 4).Value = 39
$rushing.Cells.Item(7, 5).Value = 10
$rushing.Cells.Item(7, 6).Value = 17

# ---------------------------------------------------------------------------
# Sheet "Receiving": no new rows, just cumulative stat bumps.
# ---------------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# Row 3: L.Murray
$receiving.Cells.Item(3, 3).Value = 13
$receiving.Cells.Item(3, 4).Value = 10

# Row 4: D.Freeman
$receiving.Cells.Item(4, 3).Value = 38
$receiving.Cells.Item(4, 4).Value = 31

# Row 6: M.Brown
$receiving.Cells.Item(6, 3).Value = 90
$receiving.Cells.Item(6, 4).Value = 70
$receiving.Cells.Item(6, 5).Value = 40
$receiving.Cells.Item(6, 7).Value = 15

# Row 9: J.Proche
$receiving.Cells.Item(9, 3).Value = 17
$receiving.Cells.Item(9, 4).Value = 14

# Row 10: T.Wallace
$receiving.Cells.Item(10, 3).Value = 4
$receiving.Cells.Item(10, 4).Value = 2

# Row 11: R.Bateman
$receiving.Cells.Item(11, 3).Value = 40
$receiving.Cells.Item(11, 4).Value = 29
$receiving.Cells.Item(11, 7).Value = 6
$receiving.Cells.Item(11, 8).Value = 5

# Row 13: M.Andrews
$receiving.Cells.Item(13, 3).Value = 103
$receiving.Cells.Item(13, 4).Value = 75
$receiving.Cells.Item(13, 5).Value = 29
$receiving.Cells.Item(13, 6).Value = 18
$receiving.Cells.Item(13, 7).Value = 21
$receiving.Cells.Item(13, 8).Value = 14

Write-Output "done"
